# EZPW_YR_FIN.xlsx update — add newest fiscal-year column (FY2018, period
# ending 2018-09-29) in front of the existing year columns on the single
# "EZPW" sheet, and correct a handful of prior-year cash-flow figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a new blank column at D; this pushes the existing D:K data
#    (and their cell styles) one column to the right, to E:L.
# ---------------------------------------------------------------------
$ws.Columns("D:D").Insert()

# ---------------------------------------------------------------------
# 2) The freshly inserted column D has no number formatting yet (Excel
#    seeds it from column C). Clone the formatting (date format on the
#    "Period Ending" rows, "#,##0" on the data rows) from column E, which
#    still carries the correct per-row styles after the shift above.
# ---------------------------------------------------------------------
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Populate column D with the new fiscal-year figures.
# ---------------------------------------------------------------------

# Income Statement (new "Period Ending" date: 2018-09-29)
$ws.Range("D7").Value2 = 43373
$ws.Range("D8").Value2 = 813500
$ws.Range("D9").Value2 = 330600
$ws.Range("D10").Value2 = 482900
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 11700
$ws.Range("D15").Value2 = 25500
$ws.Range("D17").Value2 = 756600
$ws.Range("D18").Value2 = 56900
$ws.Range("D20").Value2 = 28000
$ws.Range("D21").Value2 = 110400
$ws.Range("D22").Value2 = 27800
$ws.Range("D23").Value2 = 57100
$ws.Range("D24").Value2 = 13400
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 43600
$ws.Range("D27").Value2 = 44600
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = -5600
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -28000
$ws.Range("D33").Value2 = 39100
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 39100

# Balance Sheet (new "Period Ending" date: 2018-09-29)
$ws.Range("D38").Value2 = 43373
$ws.Range("D41").Value2 = 286000
$ws.Range("D42").Value2 = "NA"
$ws.Range("D43").Value2 = 271000
$ws.Range("D44").Value2 = 167000
$ws.Range("D45").Value2 = 33200
$ws.Range("D46").Value2 = 757100
$ws.Range("D47").Value2 = 52700
$ws.Range("D48").Value2 = 73600
$ws.Range("D49").Value2 = 352400
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 11000
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 1246900
$ws.Range("D57").Value2 = 57800
$ws.Range("D58").Value2 = 192700
$ws.Range("D59").Value2 = 9300
$ws.Range("D60").Value2 = 259800
$ws.Range("D61").Value2 = 226700
$ws.Range("D62").Value2 = 15700
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 498900
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 392200
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 748000
$ws.Range("D77").Value2 = 0

# Cash Flow Statement (new "Period Ending" date: 2018-09-29)
$ws.Range("D80").Value2 = 43373
$ws.Range("D81").Value2 = 39100
$ws.Range("D83").Value2 = 25500
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 88700
$ws.Range("D91").Value2 = -40500
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -134200
$ws.Range("D96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = 167600
$ws.Range("D101").Value2 = -500
$ws.Range("D102").Value2 = 121600

# ---------------------------------------------------------------------
# 4) A few prior-year (now column E) cash-flow figures were also revised
#    as part of this update, not just shifted over from the old column D.
# ---------------------------------------------------------------------
$ws.Range("E89").Value2 = 58000
$ws.Range("E91").Value2 = -25000
$ws.Range("E94").Value2 = -13400
